$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.482.45'
$ws.Range("E2").Value = '  +2.23%  '
$ws.Range("D3").Value = '1.872.25'
$ws.Range("E3").Value = '  +1.51%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.016'
$ws.Range("E4").Value = '  +0.78%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.91'
$ws.Range("E5").Value = '  +1.48%  '
$ws.Range("E6").Value = '  +0.81%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4793'
$ws.Range("E7").Value = '  +1.46%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3765'
$ws.Range("E8").Value = '  +2.81%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07382'
$ws.Range("E9").Value = '  +2.87%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9400'
$ws.Range("E10").Value = '  +1.89%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.68'
$ws.Range("E11").Value = '  +5.57%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07890'
$ws.Range("E12").Value = '  +3.81%  '
$ws.Range("D13").Value = '1.890.63'
$ws.Range("E13").Value = '  +2.62%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.440'
$ws.Range("E14").Value = '  +2.84%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.605'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '91.00'
$ws.Range("E16").Value = '  +2.93%  '
$ws.Range("E17").Value = '  +0.71%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008948'
$ws.Range("E18").Value = '  +3.58%  '
$ws.Range("E20").Value = '  +2.74%  '
$ws.Range("D21").Value = '27.514.68'
$ws.Range("E21").Value = '  +2.21%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.141'
$ws.Range("E22").Value = '  +2.58%  '
$ws.Range("E23").Value = '  +0.55%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.962'
$ws.Range("E24").Value = '  +2.45%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '154.04'
$ws.Range("E25").Value = '  +1.35%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '18.57'
$ws.Range("E26").Value = '  +2.30%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.019'
$ws.Range("E27").Value = '  +0.78%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '116.16'
$ws.Range("E28").Value = '  +1.54%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.009'
$ws.Range("E29").Value = '  +3.00%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.08929'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.328'
$ws.Range("E31").Value = '  +0.86%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.217'
$ws.Range("E32").Value = '  +4.53%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.601'
$ws.Range("E33").Value = '  +2.57%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7487'
$ws.Range("E34").Value = '  +0.40%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.691'
$ws.Range("E35").Value = '  -3.36%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02071'
$ws.Range("E36").Value = '  +6.24%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.121'
$ws.Range("E37").Value = '  +3.03%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05296'
$ws.Range("E38").Value = '  +0.77%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.998'
$ws.Range("E39").Value = '  +1.15%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5376'
$ws.Range("E40").Value = '  +3.29%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.083'
$ws.Range("E41").Value = '  +2.83%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1530'
$ws.Range("E42").Value = '  +1.29%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.433'
$ws.Range("E43").Value = '  +3.12%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.4843'
$ws.Range("E44").Value = '  +3.27%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.59'
$ws.Range("E45").Value = '  +0.69%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.015'
$ws.Range("E46").Value = '  +0.88%  '
$ws.Range("E47").Value = '  +4.07%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '103.18'
$ws.Range("E48").Value = '  +1.42%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '67.09'
$ws.Range("E49").Value = '  +2.54%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06112'
$ws.Range("E50").Value = '  +1.44%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.9006'
$ws.Range("E51").Value = '  +1.88%  '
